$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 value from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8 to match the resulting sheetView selection
$ws.Range("E8").Select()
